$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 46079, 12.406, 0.267),
    @(3, 46079.01041666666, 2.231, 0.748),
    @(4, 46079.02083333334, 3.8, 0.217),
    @(5, 46079.03125, 3.092, 0.457),
    @(6, 46079.04166666666, 17.412, 0.018),
    @(7, 46079.05208333334, 24.545, 0),
    @(8, 46079.0625, 30.145, 0),
    @(9, 46079.07291666666, 14.545, 0.104),
    @(10, 46079.08333333334, 6.692, 0.099),
    @(11, 46079.09375, 12.518, 0),
    @(12, 46079.10416666666, 1.762, 0.004),
    @(13, 46079.11458333334, 10.892, 0),
    @(14, 46079.125, 11.2, 0),
    @(15, 46079.13541666666, 11.882, 0),
    @(16, 46079.14583333334, 3.459, 1.39),
    @(17, 46079.15625, 0, 4.762),
    @(18, 46079.16666666666, 0, 2.983),
    @(19, 46079.17708333334, 0, 2.74),
    @(20, 46079.1875, 15.189, 0.124),
    @(21, 46079.19791666666, 11.738, 0),
    @(22, 46079.20833333334, 23.345, 0),
    @(23, 46079.21875, 23.495, 0),
    @(24, 46079.22916666666, 15.371, 0),
    @(25, 46079.23958333334, 22.128, 0),
    @(26, 46079.25, 3.152, 0.331),
    @(27, 46079.26041666666, 31.594, 0),
    @(28, 46079.27083333334, 48.481, 0),
    @(29, 46079.28125, 40.719, 0),
    @(30, 46079.29166666666, 46.416, 0),
    @(31, 46079.30208333334, 6.198, 0.407),
    @(32, 46079.3125, 0, 1.22),
    @(33, 46079.32291666666, 0, 0.141),
    @(34, 46079.33333333334, 0.186, 0.975),
    @(35, 46079.34375, 0, 4.31),
    @(36, 46079.35416666666, 0, 1.118),
    @(37, 46079.36458333334, 0, 0.464),
    @(38, 46079.375, 7.07, 0.533),
    @(39, 46079.38541666666, 17.821, 0),
    @(40, 46079.39583333334, 4.435, 0.136),
    @(41, 46079.40625, 0, 0.508),
    @(42, 46079.41666666666, 0, 0),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Newly created rows (32-42) need the same timestamp number format as the
# existing rows in column A ("YYYY-MM-DD HH:MM:SS"); reapplying it across the
# whole data range keeps a single consistent style for all timestamp cells.
$ws.Range("A2:A42").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Done updating rows"